$wb = $excel.ActiveWorkbook

# Old / new values used in this "Generate Report for Handoff" run.
$oldGuid = "677ca8d3-5ad7-4d75-a45d-6451a5f3def1"
$newGuid = "3e741122-46f3-4ed5-834c-a68700b9b525"

$oldHash = "37f061e3a6571e3e0557f87d920896e51aa1794e"
$newHash = "b79ccbd1786e20b0c71f1aaae5b7a29ba6a0cb25"

$oldZhDatetime = "2016-03-11 03:28:25"
$newZhDatetime = "2016-03-11 03:28:58"

$oldDeDatetime = "2016-03-11 03:28:32"
$newDeDatetime = "2016-03-11 03:29:05"

$newMdName = $newGuid + ".md"
$newZhXlf = $newGuid + "." + $newHash + ".zh-cn.xlf"
$newDeXlf = $newGuid + "." + $newHash + ".de-de.xlf"

# 1) Update the actual cell contents (this keeps the shared-string table
#    consistent - once every usage of an old string is gone it drops out).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace($oldGuid, $newGuid)
    $ws.Cells.Replace($oldHash, $newHash)
    $ws.Cells.Replace($oldZhDatetime, $newZhDatetime)
    $ws.Cells.Replace($oldDeDatetime, $newDeDatetime)
}

# 2) Update the hyperlinks' display text to match (Replace() above only
#    touches cell values, not the hyperlink "display" text).
foreach ($ws in $wb.Worksheets) {
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$2') {
            $h.TextToDisplay = $newMdName
        } elseif ($addr -eq '$C$2') {
            if ($ws.Name -eq "zh-cn") {
                $h.TextToDisplay = $newZhXlf
            } elseif ($ws.Name -eq "de-de") {
                $h.TextToDisplay = $newDeXlf
            }
        }
    }
}
